$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1939
$ws.Range("I2").Value = 5196
$ws.Range("J2").Value = 20801
$ws.Range("K2").Value = 96
$ws.Range("L2").Value = 5558
$ws.Range("M2").Value = 314
$ws.Range("N2").Value = 3581
$ws.Range("O2").Value = 10
$ws.Range("P2").Value = 77
$ws.Range("Q2").Value = 33
$ws.Range("R2").Value = 280
$ws.Range("S2").Value = 2224
$ws.Range("T2").Value = 3610
$ws.Range("U2").Value = 276
$ws.Range("V2").Value = 32433
$ws.Range("W2").Value = 15
$ws.Range("X2").Value = 32310
$ws.Range("Y2").Value = 41
$ws.Range("Z2").Value = 496
$ws.Range("AA2").Value = 218
